$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay plain text (workbook stores them as
# inline strings, e.g. "1.00", "0.0439") instead of being auto-converted
# to numbers by Excel's smart-typing on assignment.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '69.139.10'
$ws.Range("E2").Value = '  +5.12%  '
$ws.Range("D3").Value = '3.542.18'
$ws.Range("E3").Value = '  +3.93%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '590.50'
$ws.Range("E5").Value = '  +4.90%  '
$ws.Range("D6").Value = '193.39'
$ws.Range("E6").Value = '  +9.31%  '
$ws.Range("D7").Value = '0.639'
$ws.Range("E7").Value = '  +0.91%  '
$ws.Range("D8").Value = '3.536.55'
$ws.Range("E8").Value = '  +4.00%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").Value = '  +4.25%  '
$ws.Range("D11").Value = '0.660'
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("D12").Value = '58.78'
$ws.Range("E12").Value = '  +9.30%  '
$ws.Range("E13").Value = '  +4.55%  '
$ws.Range("D14").Value = '9.63'
$ws.Range("E14").Value = '  +4.03%  '
$ws.Range("D15").Value = '4.086.91'
$ws.Range("E15").Value = '  +3.49%  '
$ws.Range("D16").Value = '19.13'
$ws.Range("E16").Value = '  +4.11%  '
$ws.Range("D17").Value = '3.540.79'
$ws.Range("E17").Value = '  +3.24%  '
$ws.Range("D18").Value = '69.143.20'
$ws.Range("E18").Value = '  +5.29%  '
$ws.Range("D19").Value = '12.34'
$ws.Range("E19").Value = '  +3.73%  '
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("D22").Value = '492.28'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").Value = '5.60'
$ws.Range("E23").Value = '  +12.99%  '
$ws.Range("D24").Value = '17.30'
$ws.Range("E24").Value = '  +21.20%  '
$ws.Range("D25").Value = '4.45'
$ws.Range("E25").Value = '  +7.54%  '
$ws.Range("D26").Value = '90.76'
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("E27").Value = '  +4.13%  '
$ws.Range("D28").Value = '11.20'
$ws.Range("E28").Value = '  +4.34%  '
$ws.Range("D29").Value = '9.20'
$ws.Range("E29").Value = '  +4.93%  '
$ws.Range("D30").Value = '31.83'
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").Value = '7.47'
$ws.Range("E31").Value = '  +13.51%  '
$ws.Range("D32").Value = '614.07'
$ws.Range("E32").Value = '  +6.59%  '
$ws.Range("D33").Value = '11.99'
$ws.Range("E33").Value = '  +3.87%  '
$ws.Range("D34").Value = '65.15'
$ws.Range("E34").Value = '  +4.22%  '
$ws.Range("E35").Value = '  +4.90%  '
$ws.Range("E36").Value = '  +5.00%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '37.77'
$ws.Range("E38").Value = '  +4.67%  '
$ws.Range("E39").Value = '  +5.19%  '
$ws.Range("E40").Value = '  +6.37%  '
$ws.Range("D41").Value = '3.56'
$ws.Range("E41").Value = '  -1.48%  '
$ws.Range("D42").Value = '3.283.46'
$ws.Range("E42").Value = '  +5.02%  '
$ws.Range("D43").Value = '2.97'
$ws.Range("E43").Value = '  +5.96%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.66'
$ws.Range("E44").Value = '  +8.88%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0439'
$ws.Range("E45").Value = '  +4.84%  '
$ws.Range("D46").Value = '3.35'
$ws.Range("E46").Value = '  +5.94%  '
$ws.Range("D47").Value = '2.78'
$ws.Range("E47").Value = '  +19.28%  '
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("D49").Value = '9.05'
$ws.Range("E49").Value = '  +6.51%  '
$ws.Range("D50").Value = '1.00'
$ws.Range("D51").Value = '141.78'
$ws.Range("E51").Value = '  +0.91%  '
